$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "Sunbury"
$ws.Range("B2").Value = "Aldente Deli - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429"
$ws.Range("C2").Value = "15:45-16:23  5/2/2021"
$ws.Range("D2").Value = "Case attended venue"
$ws.Range("E2").Value = "new"

# --- Row 3 ---
$ws.Range("A3").Value = "Sunbury"
$ws.Range("B3").Value = "Asian Star - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429"
$ws.Range("C3").Value = "15:57-16:30  5/2/2021"
$ws.Range("D3").Value = "Case attended venue"
$ws.Range("E3").Value = "new"

# --- Row 4 ---
$ws.Range("A4").Value = "Sunbury"
$ws.Range("B4").Value = "Bakers Delight - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429"
$ws.Range("C4").Value = "15:40-16:15  5/2/2021"
$ws.Range("D4").Value = "Case attended venue"
$ws.Range("E4").Value = "new"

# --- Row 5 ---
$ws.Range("A5").Value = "Sunbury"
$ws.Range("B5").Value = "Cellarbrations  34 Batman Avenue  Sunbury VIC 3429"
$ws.Range("C5").Value = "17:44-18:19  7/2/2021"
$ws.Range("D5").Value = "Case attended venue"
$ws.Range("E5").Value = "new"

# --- Row 6 ---
$ws.Range("A6").Value = "Sunbury"
$ws.Range("B6").Value = "Cellarbrations  34 Batman Avenue  Sunbury VIC 3429"
$ws.Range("C6").Value = "18:17-19:02  6/2/2021"
$ws.Range("D6").Value = "Case attended venue"
$ws.Range("E6").Value = "new"

# --- Row 7 ---
$ws.Range("A7").Value = "Sunbury"
$ws.Range("B7").Value = "PJ's Pet Warehouse  Shop 2, 104 Horne Street  Sunbury VIC 3429"
$ws.Range("C7").Value = "15:37-16:10  5/2/2021"
$ws.Range("D7").Value = "Case attended venue"
$ws.Range("E7").Value = "new"

# --- Row 8 ---
$ws.Range("A8").Value = "Sunbury"
$ws.Range("B8").Value = "Sunny Life Massage - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429"
$ws.Range("C8").Value = "16:30-18:30  6/2/2021"
$ws.Range("D8").Value = "Case attended venue"
$ws.Range("E8").Value = "new"

# --- Row 9 ---
$ws.Range("A9").Value = "Sunbury"
$ws.Range("B9").Value = "Sushi Sushi - Sunbury Square Shopping Centre  2-28 Evans Street  Sunbury VIC 3429"
$ws.Range("C9").Value = "15:53-16:28  5/2/2021"
$ws.Range("D9").Value = "Case attended venue"
$ws.Range("E9").Value = "new"

# Resize columns A and B to fit the new (longer) content, mirroring the
# AutoFit Excel performs when the data changes.
$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()

# Match the page setup (paper size / orientation) recorded for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection ends up on the new last data cell in column C, matching the
# author's on-screen state when the file was saved.
$ws.Range("C2:C9").Select()
